# Apply the cryptos list data refresh (new prices / 1h volume%)
# plus the 4 rank-position swaps captured in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'30.240.20"
$ws.Range('E2').Value = '  +1.01%  '

$ws.Range('D3').Value = "'1.673.14"
$ws.Range('E3').Value = '  +2.65%  '

$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = '  -0.30%  '

$ws.Range('D5').Value = "'218.19"
$ws.Range('E5').Value = '  +1.73%  '

$ws.Range('D6').Value = "'0.523"
$ws.Range('E6').Value = '  +0.09%  '

$ws.Range('D7').Value = "'0.994"
$ws.Range('E7').Value = '  -0.61%  '

$ws.Range('D8').Value = "'29.47"
$ws.Range('E8').Value = '  +0.47%  '

$ws.Range('D9').Value = "'0.266"
$ws.Range('E9').Value = '  +2.37%  '

$ws.Range('D10').Value = "'0.0619"
$ws.Range('E10').Value = '  +0.81%  '

$ws.Range('D11').Value = "'0.0904"
$ws.Range('E11').Value = '  -1.44%  '

$ws.Range('D12').Value = "'1.914.88"
$ws.Range('E12').Value = '  +2.71%  '

$ws.Range('B13').Value = 'Chainlink'
$ws.Range('C13').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D13').Value = "'10.44"
$ws.Range('E13').Value = '  +15.43%  '

$ws.Range('D14').Value = "'0.617"
$ws.Range('E14').Value = '  +7.82%  '

$ws.Range('B15').Value = 'WrappedEther'
$ws.Range('C15').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D15').Value = "'1.662.82"
$ws.Range('E15').Value = '  +1.76%  '

$ws.Range('D16').Value = "'4.00"
$ws.Range('E16').Value = '  +1.97%  '

$ws.Range('D17').Value = "'30.291.55"
$ws.Range('E17').Value = '  +1.06%  '

$ws.Range('D18').Value = "'65.63"
$ws.Range('E18').Value = '  +1.37%  '

$ws.Range('D19').Value = "'245.59"
$ws.Range('E19').Value = '  -0.70%  '

$ws.Range('D20').Value = "'0.0₃0716"
$ws.Range('E20').Value = '  +1.27%  '

$ws.Range('D21').Value = "'0.997"
$ws.Range('E21').Value = '  -0.35%  '

$ws.Range('D22').Value = "'4.30"
$ws.Range('E22').Value = '  +3.97%  '

$ws.Range('D23').Value = "'10.06"
$ws.Range('E23').Value = '  +4.38%  '

$ws.Range('D24').Value = "'2.19"
$ws.Range('E24').Value = '  +3.37%  '

$ws.Range('D25').Value = "'158.33"
$ws.Range('E25').Value = '  -0.05%  '

$ws.Range('D26').Value = "'15.92"
$ws.Range('E26').Value = '  +0.90%  '

$ws.Range('D27').Value = "'0.111"
$ws.Range('E27').Value = '  -0.17%  '

$ws.Range('D28').Value = "'6.73"
$ws.Range('E28').Value = '  +1.43%  '

$ws.Range('D29').Value = "'0.998"
$ws.Range('E29').Value = '  -0.28%  '

$ws.Range('D30').Value = "'0.0500"
$ws.Range('E30').Value = '  +1.62%  '

$ws.Range('B31').Value = 'Filecoin'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D31').Value = "'3.47"
$ws.Range('E31').Value = '  +3.19%  '

$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').Value = "'1.12"
$ws.Range('E32').Value = '  -0.38%  '

$ws.Range('D33').Value = "'3.26"
$ws.Range('E33').Value = '  +1.71%  '

$ws.Range('D34').Value = "'1.460.24"
$ws.Range('E34').Value = '  +2.09%  '

$ws.Range('D35').Value = "'1.73"
$ws.Range('E35').Value = '  +5.18%  '

$ws.Range('D36').Value = "'1.03"
$ws.Range('E36').Value = '  -1.06%  '

$ws.Range('D37').Value = "'0.0178"
$ws.Range('E37').Value = '  +4.16%  '

$ws.Range('D38').Value = "'80.16"
$ws.Range('E38').Value = '  +14.45%  '

$ws.Range('D39').Value = "'0.586"
$ws.Range('E39').Value = '  +4.59%  '

$ws.Range('B40').Value = 'MXToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D40').Value = "'2.69"
$ws.Range('E40').Value = '  -7.42%  '

$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').Value = "'2.30"
$ws.Range('E41').Value = '  +0.28%  '

$ws.Range('D42').Value = "'0.857"
$ws.Range('E42').Value = '  +2.64%  '

$ws.Range('D43').Value = "'0.0503"
$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('D44').Value = "'1.99"
$ws.Range('E44').Value = '  +0.86%  '

$ws.Range('B45').Value = 'BitcoinSV'
$ws.Range('C45').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D45').Value = "'53.10"
$ws.Range('E45').Value = '  -3.01%  '

$ws.Range('D46').Value = "'1.02"
$ws.Range('E46').Value = '  -2.23%  '

$ws.Range('B47').Value = 'PaxDollar'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D47').Value = "'0.995"
$ws.Range('E47').Value = '  -0.48%  '

$ws.Range('B48').Value = 'RocketPoolETH'
$ws.Range('C48').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D48').Value = "'1.810.07"
$ws.Range('E48').Value = '  +2.09%  '

$ws.Range('D49').Value = "'5.42"
$ws.Range('E49').Value = '  -0.60%  '

$ws.Range('D50').Value = "'95.45"
$ws.Range('E50').Value = '  +6.65%  '

$ws.Range('D51').Value = "'0.0₆0110"
$ws.Range('E51').Value = '  +2.28%  '
